# Fruta / hortaliza, semanal
# Insert two new weekly observation rows into the daily-logic subset sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows starting at row 359; everything from the old
# row 359 downward shifts down by two rows (old 359 -> 361, old 360 -> 362, ...).
$ws.Rows("359:360").Insert()

# --- New row 359 --------------------------------------------------------
$ws.Range("A359").Value = 9
$ws.Range("B359").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C359").Value = "Metropolitana"
$ws.Range("D359").Value = 44798
$ws.Range("E359").Value = 13
$ws.Range("F359").Value = 100112039
$ws.Range("G359").Value = "Ciboulette"
$ws.Range("H359").Value = "Sin especificar"
$ws.Range("I359").Value = "Primera"
$ws.Range("J359").Value = 250
$ws.Range("K359").Value = 2000
$ws.Range("L359").Value = 2000
$ws.Range("M359").Value = 2000
$ws.Range("N359").Value = "`$/docena de atados"
$ws.Range("O359").Value = "Región Metropolitana"
$ws.Range("P359").Value = 667
$ws.Range("Q359").Value = 3
$ws.Range("R359").Value = "Hortaliza"

# --- New row 360 --------------------------------------------------------
$ws.Range("A360").Value = 9
$ws.Range("B360").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C360").Value = "Metropolitana"
$ws.Range("D360").Value = 44798
$ws.Range("E360").Value = 13
$ws.Range("F360").Value = 100112039
$ws.Range("G360").Value = "Ciboulette"
$ws.Range("H360").Value = "Sin especificar"
$ws.Range("I360").Value = "Segunda"
$ws.Range("J360").Value = 160
$ws.Range("K360").Value = 1500
$ws.Range("L360").Value = 1500
$ws.Range("M360").Value = 1500
$ws.Range("N360").Value = "`$/docena de atados"
$ws.Range("O360").Value = "Región Metropolitana"
$ws.Range("P360").Value = 500
$ws.Range("Q360").Value = 3
$ws.Range("R360").Value = "Hortaliza"
